{"js": "// Office.js (Word JavaScript API) script \u2014 applies the report.docx edits\n// described in the commit diff: softmax/LayerNorm wording updates, the\n// resource-utilisation table row for \"Layer Normalisation\", and the\n// Limitations/Conclusion paragraph rewrites.\n//\n// Short, unique anchor substrings are used with body.search() and then\n// expanded to the enclosing paragraph via `.paragraphs.getFirst()` before\n// the whole paragraph text is overwritten \u2014 this keeps each search string\n// well under any host search-length limits and lets one call collapse a\n// paragraph made of several runs (e.g. the \"7.2 Remaining\" bullet) down to\n// the single new run the diff expects.\n\nasync function replaceParagraphByAnchor(context, anchor, newText) {\n  const body = context.document.body;\n  const results = body.search(anchor, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  const para = results.items[0].paragraphs.getFirst();\n  para.insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Summary bullet: \"division-free softmax normalisation\" -> \"division-free softmax and LayerNorm\"\nawait replaceParagraphByAnchor(\n  context,\n  \"Key design features include Q8.8 fixed-p\",\n  \"Key design features include Q8.8 fixed-point arithmetic, systolic array compute, BRAM-backed weight and KV-cache storage, division-free softmax and LayerNorm via reciprocal LUT with Newton-Raphson refinement, and packed array ports for full simulation visibility.\"\n);\n\n// 2) 3.3.4 Layer Normalization paragraph \u2014 expanded explanation of the\n// division-free right-shift / rsqrt LUT / Newton-Raphson design.\nawait replaceParagraphByAnchor(\n  context,\n  \"LayerNorm computes the mean and variance\",\n  \"LayerNorm computes the mean and variance of the input vector, then normalizes each element through three FSM stages: mean computation via sequential accumulation and arithmetic right-shift by log\u2082(VEC_LEN), variance computation using centered differences and the same right-shift, and element-wise normalization with learnable gamma and beta parameters. The reciprocal square root (1/\u221a(variance + \u03b5)) uses a 32-entry LUT indexed by CLZ-normalised mantissa bits, followed by one Newton-Raphson iteration \u2014 the same architectural pattern proven in the softmax normalisation path. No runtime division operators remain; all \u00f7N operations use arithmetic right-shift since VEC_LEN is a power of 2.\"\n);\n\n// 3) Resource-utilisation table, \"Layer Normalisation\" row: DSP48 / FFs / LUTs\n//    columns. Addressed by row/column index (not text search) because \"~500\"\n//    also appears, unrelated, in the \"4x4 Systolic Array\" row.\n{\n  const tables = context.document.body.tables;\n  tables.load(\"items\");\n  await context.sync();\n\n  const table = tables.items[1];\n  table.load(\"rowCount\");\n  await context.sync();\n\n  const rowCount = table.rowCount;\n  const firstCells = [];\n  for (let i = 0; i < rowCount; i++) {\n    const c = table.getCell(i, 0);\n    c.body.load(\"text\");\n    firstCells.push(c);\n  }\n  await context.sync();\n\n  let targetRow = -1;\n  for (let i = 0; i < rowCount; i++) {\n    if (firstCells[i].body.text.trim() === \"Layer Normalisation\") {\n      targetRow = i;\n      break;\n    }\n  }\n\n  if (targetRow >= 0) {\n    table.getCell(targetRow, 1).getRange().insertText(\"4\u20136\", Word.InsertLocation.replace);\n    table.getCell(targetRow, 2).getRange().insertText(\"~400\", Word.InsertLocation.replace);\n    table.getCell(targetRow, 3).getRange().insertText(\"~600\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 4) \"7.1 Resolved\" paragraph \u2014 append the new division-free LayerNorm clause.\nawait replaceParagraphByAnchor(\n  context,\n  \"BRAM-backed weights, division-free softm\",\n  \"BRAM-backed weights, division-free softmax, packed array ports, streaming weight architecture (99.4% FF reduction), division-free LayerNorm (arithmetic right-shift for \u00f7N, 32-entry rsqrt LUT + Newton-Raphson for 1/\u221avar).\"\n);\n\n// 5) \"7.2 Remaining\" paragraph \u2014 the old wording was split across five runs\n// (one styled \"/\" run for the division operator). The anchor text is unique\n// to this paragraph, so expanding to the enclosing paragraph and overwriting\n// collapses it to the single new run the diff expects.\nawait replaceParagraphByAnchor(\n  context,\n  \"LayerNorm still uses\",\n  \"Only ReLU activation supported (not GELU/SiLU). Systolic array instantiated but not connected to datapath. Single-token processing with no batching.\"\n);\n\n// 6) \"7.3 Planned\" paragraph.\nawait replaceParagraphByAnchor(\n  context,\n  \"LayerNorm right-shift + reciprocal-LUT, \",\n  \"Systolic-tiled projections, parallel softmax with N_HEADS instances, GELU/SiLU activation via PWL approximation, multi-layer stacking, multi-device distribution, AXI-Lite control/status interface.\"\n);\n\n// 7) \"8. Conclusion\" paragraph \u2014 insert the division-free datapath clause.\nawait replaceParagraphByAnchor(\n  context,\n  \"This project demonstrates a complete, sy\",\n  \"This project demonstrates a complete, synthesizable transformer decoder block in SystemVerilog with two architectural variants: a high-throughput register-bridge design and a minimum-area streaming design. The 17-module hierarchy mirrors the transformer\u2019s conceptual structure while addressing practical synthesis concerns: no runtime division operators remain in the compute datapath (softmax and LayerNorm both use LUT + Newton-Raphson, mean/variance use arithmetic right-shift). All 83 verification tests pass across both variants. The streaming architecture achieves 99.4% register reduction while maintaining functional equivalence, making the design feasible on the smallest FPGA targets.\"\n);\n", "ps1": "# Word COM interop script \u2014 applies the report.docx edits described in the\n# commit diff: softmax/LayerNorm wording updates, the resource-utilisation\n# table row for \"Layer Normalisation\", and the Limitations/Conclusion\n# paragraph rewrites.\n#\n# NOTE: each Find/Replace is inlined (no helper function with named params)\n# because this host's PowerShell-subset does not bind `-Name value` args to\n# custom function parameters.\n\n$d = $word.ActiveDocument\n\n# 1) Summary bullet: \"division-free softmax normalisation\" -> \"division-free softmax and LayerNorm\"\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"Key design features include Q8.8 fixed-point arithmetic, systolic array compute, BRAM-backed weight and KV-cache storage, division-free softmax normalisation via reciprocal LUT with Newton-Raphson refinement, and packed array ports for full simulation visibility.\"\n$find.Replacement.Text = \"Key design features include Q8.8 fixed-point arithmetic, systolic array compute, BRAM-backed weight and KV-cache storage, division-free softmax and LayerNorm via reciprocal LUT with Newton-Raphson refinement, and packed array ports for full simulation visibility.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2) 3.3.4 Layer Normalization paragraph \u2014 expanded explanation of the\n# division-free right-shift / rsqrt LUT / Newton-Raphson design.\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"LayerNorm computes the mean and variance of the input vector, then normalizes each element through three FSM stages: mean computation via sequential accumulation, variance computation using centered differences, and element-wise normalization with learnable gamma and beta parameters. The reciprocal square root uses a 4-entry LUT.\"\n$find.Replacement.Text = \"LayerNorm computes the mean and variance of the input vector, then normalizes each element through three FSM stages: mean computation via sequential accumulation and arithmetic right-shift by log\u2082(VEC_LEN), variance computation using centered differences and the same right-shift, and element-wise normalization with learnable gamma and beta parameters. The reciprocal square root (1/\u221a(variance + \u03b5)) uses a 32-entry LUT indexed by CLZ-normalised mantissa bits, followed by one Newton-Raphson iteration \u2014 the same architectural pattern proven in the softmax normalisation path. No runtime division operators remain; all \u00f7N operations use arithmetic right-shift since VEC_LEN is a power of 2.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 3) Resource-utilisation table, \"Layer Normalisation\" row: DSP48 / FFs / LUTs\n#    columns. Addressed by row/column index (not text search) because \"~500\"\n#    also appears, unrelated, in the \"4x4 Systolic Array\" row.\n$table = $d.Tables.Item(2)\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    $label = $table.Cell($r, 1).Range.Text.TrimEnd([char]13, [char]7)\n    if ($label -eq \"Layer Normalisation\") {\n        $c1 = $table.Cell($r, 2).Range\n        $c1.End = $c1.End - 1\n        $c1.Text = \"4\u20136\"\n\n        $c2 = $table.Cell($r, 3).Range\n        $c2.End = $c2.End - 1\n        $c2.Text = \"~400\"\n\n        $c3 = $table.Cell($r, 4).Range\n        $c3.End = $c3.End - 1\n        $c3.Text = \"~600\"\n    }\n}\n\n# 4) \"7.1 Resolved\" paragraph \u2014 append the new division-free LayerNorm clause.\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"BRAM-backed weights, division-free softmax, packed array ports, streaming weight architecture (99.4% FF reduction).\"\n$find.Replacement.Text = \"BRAM-backed weights, division-free softmax, packed array ports, streaming weight architecture (99.4% FF reduction), division-free LayerNorm (arithmetic right-shift for \u00f7N, 32-entry rsqrt LUT + Newton-Raphson for 1/\u221avar).\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 5) \"7.2 Remaining\" paragraph \u2014 the old wording was split across five runs\n# (one styled \"/\" run for the division operator). Find the anchor text,\n# expand to the whole paragraph, and overwrite with the single new run.\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"LayerNorm still uses\"\n$find.Execute() | Out-Null\n$rng.Expand(4) | Out-Null      # wdParagraph\n$rng.MoveEnd(1, -1) | Out-Null # exclude the trailing paragraph mark\n$rng.Text = \"Only ReLU activation supported (not GELU/SiLU). Systolic array instantiated but not connected to datapath. Single-token processing with no batching.\"\n\n# 6) \"7.3 Planned\" paragraph.\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"LayerNorm right-shift + reciprocal-LUT, systolic-tiled projections, parallel softmax, GELU/SiLU PWL, multi-layer stacking, multi-device distribution, AXI-Lite interface.\"\n$find.Replacement.Text = \"Systolic-tiled projections, parallel softmax with N_HEADS instances, GELU/SiLU activation via PWL approximation, multi-layer stacking, multi-device distribution, AXI-Lite control/status interface.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 7) \"8. Conclusion\" paragraph \u2014 insert the division-free datapath clause.\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"This project demonstrates a complete, synthesizable transformer decoder block in SystemVerilog with two architectural variants: a high-throughput register-bridge design and a minimum-area streaming design. The 17-module hierarchy mirrors the transformer\u2019s conceptual structure while addressing practical synthesis concerns. All 83 verification tests pass across both variants. The streaming architecture achieves 99.4% register reduction while maintaining functional equivalence, making the design feasible on the smallest FPGA targets.\"\n$find.Replacement.Text = \"This project demonstrates a complete, synthesizable transformer decoder block in SystemVerilog with two architectural variants: a high-throughput register-bridge design and a minimum-area streaming design. The 17-module hierarchy mirrors the transformer\u2019s conceptual structure while addressing practical synthesis concerns: no runtime division operators remain in the compute datapath (softmax and LayerNorm both use LUT + Newton-Raphson, mean/variance use arithmetic right-shift). All 83 verification tests pass across both variants. The streaming architecture achieves 99.4% register reduction while maintaining functional equivalence, making the design feasible on the smallest FPGA targets.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
